# Updates cryptos list figures (price + 1h volume change) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '31.044.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.02%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.959.86'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.19%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.56%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4873'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2955'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.80%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06968'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.07%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.54'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.87%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '108.31'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.951.10'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.62%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07801'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.65%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.503'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.55%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7026'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.97%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '282.10'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.68%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '31.048.82'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.34'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.40%  '

$ws.Range("B19").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C19").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.250.88'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.69%  '

$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007794'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.27%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.0000'
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.548'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.19%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.40%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.540'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.892'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.24%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.95'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.29%  '

$ws.Range("E27").Value = '  -0.28%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.197'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.06%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1052'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.40%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.388'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.76%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.652'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.85%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.576'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.77%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.485'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.19%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04932'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.39%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7587'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.69%  '

$ws.Range("E36").Value = '  -0.54%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.732'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.34%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02021'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.13%  '

$ws.Range("E39").Value = '  -0.73%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.597'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '78.35'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +11.84%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.131'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.27%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9028'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.91%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '109.97'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4466'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.150'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.98%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.000'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.014.69'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.388'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.21%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1260'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.35%  '

$ws.Range("E51").Value = '  -0.01%  '
